$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.768.41"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.79"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.24"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.58"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.997"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.664.36"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.306.60"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.749.39"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.11"
$ws.Range("E21").Value = "  -9.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.63"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.50"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.35"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  +12.44%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.53"
$ws.Range("E30").Value = "  +2.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.42"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.93"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0872"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.75"
$ws.Range("E34").Value = "  +6.24%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.83"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.67"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.60"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.81"
$ws.Range("E42").Value = "  +9.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.24"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.233"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.82"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.667.79"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.71"
$ws.Range("E49").Value = "  -5.67%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.25"
$ws.Range("E51").Value = "  +0.34%  "
